$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# --- "Estatistica" block (rows 3-7): the divisor year moved from row 3 to row 4 ---
# Row 3 standalone formulas
$ws.Range("D3").Formula = '=B3*$C$4'
$ws.Range("H3").Formula = '=F3*$G$4'
$ws.Range("M3").Formula = '=K3*$L$4'
$ws.Range("Q3").Formula = '=P3/$P$4'

# Rows 4:7 shared formulas (fill whole range so Excel keeps them as one shared group)
$ws.Range("D4:D7").Formula = '=B4*$C$4'
$ws.Range("H4:H7").Formula = '=F4*$G$4'
$ws.Range("M4:M7").Formula = '=K4*$L$4'

# Q4 is now the base year (P4/P4 = 1) and stands alone
$ws.Range("Q4").Formula = '=P4/$P$4'

# Q5:Q7 become their own shared group referencing the new base ($P$4)
$ws.Range("Q5:Q7").Formula = '=P5/$P$4'

# --- New "retornos" column R: year-over-year ratio of the Q-column index ---
$ws.Range("R5").Formula = '=Q5/Q4'
$ws.Range("R5").ClearFormats()

$ws.Range("R6:R7").Formula = '=Q6/Q5'
$ws.Range("R6:R7").ClearFormats()

# --- New row 9: blank Q9 (same style as the Q column above) + PRODUCT in R9 ---
$ws.Range("Q9").NumberFormat = "0.000000"

$ws.Range("R9").Formula = '=PRODUCT(R5:R7)'
$ws.Range("R9").NumberFormat = "0.0000"

# Match the author's final selection
$ws.Range("R9").Select()
